$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.050.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.10%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.06"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.94%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.58"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.329"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0694"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0984"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.115.09"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.850.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.36"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.677"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.70"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.018.44"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.02"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0793"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.07"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.14"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.81"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.34%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.23"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.85"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +22.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.88"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.65"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.30%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.23%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0556"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.98"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.98"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.31%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +24.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.98"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +10.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.760"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +8.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.25"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.91%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.08"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +11.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "91.27"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0202"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.346.20"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.83"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.32"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.52%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Gas"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +82.30%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.39"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.03%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.76"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.80%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.37"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.68%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Kaspa"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0532"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.026.03"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.45"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +15.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0677"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.24%  "
